$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("A1").Value = "house_id"
$ws.Range("B1").Value = "stories"
$ws.Range("C1").Value = "type"
$ws.Range("D1").Value = "Bedrooms"
$ws.Range("E1").Value = "Bathrooms"
$ws.Range("F1").Value = "parking"
$ws.Range("G1").Value = "utilities"
$ws.Range("H1").Value = "Address"
$ws.Range("I1").Value = "Postal Code"

# Row 2
$ws.Range("A2").Value = 1234
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = "Townhouse"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = "yes"
$ws.Range("G2").Value = "landlord"
$ws.Range("H2").Value = "1111 Market Street, Philadelphia, PA"
$ws.Range("I2").Value = 19103

# Row 3
$ws.Range("A3").Value = 5678
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "Condo"
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = "yes"
$ws.Range("G3").Value = "tenant"
$ws.Range("H3").Value = "2222 Aldine Street, Philadelphia, PA"
$ws.Range("I3").Value = 19136

# Row 4
$ws.Range("A4").Value = 9102
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "Duplex"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "no"
$ws.Range("G4").Value = "landlord"
$ws.Range("H4").Value = "3333 Cabell Road, Philadelphia, PA"
$ws.Range("I4").Value = 19154

# Column widths (closest achievable values under this runtime's width model)
$ws.Range("D1:E1").EntireColumn.ColumnWidth = 9.75
$ws.Range("H1").EntireColumn.ColumnWidth = 31.25

# Selection / view
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("I2").Select()
